# =====================================================================
# R051-R053 Consumer Price Index (.recapitulation)
# Adds three new reference rows (271-273) to the MathDIY fundamentals
# sheet, covering: CPI for a single item, CPI for multiple items, and
# the Harmonized Index of Consumer Prices (HICP).
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---- Row 271 ----
$ws.Range("A271").Value = 'R051'
$cell = $ws.Range("B271")
$cell.Value = 'CPI2 ÷ CPI1 = p2 ÷ p1'
$run = $cell.Characters(1, 3)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(4, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(5, 6)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(11, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(12, 4)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(16, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(17, 4)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(21, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$ws.Range("C271").Value = 'Calculating the CPI for a single item'
$cell = $ws.Range("D271")
$cell.Value = '          market basket of desired JJJJ
CPI = —————————————— x 100
          market basket of base JJJJ 

           p2 of item, given period (updated cost)
CPI = ———————————————— x 100
           p1 of item, initial period (base cost
'
$run = $cell.Characters(1, 40)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(41, 27)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(68, 38)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(106, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(107, 12)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(119, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(120, 38)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(158, 29)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(187, 12)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(199, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(200, 36)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$ws.Range("E271").Value = 'Excerpt retrieved from:
en.wikipedia.org; 
title id 1062042105, Consumer Price Index 
This source was last edited on 25 Dec 2021, at 21:28 (UTC). 
Source above is available under the 
Creative Commons Attribution-ShareAlike License; additional terms may apply. By using this source, you agree to the Terms of Use and Privacy Policy. Wikipedia® is a registered trademark of the Wikimedia Foundation, Inc., a non-profit organization. '
$ws.Range("F271").Value = '.recapitulation'

# ---- Row 272 ----
$ws.Range("A272").Value = 'R052'
$cell = $ws.Range("B272")
$cell.Value = '             n 
           ∑         CPI x weight
             i=1
CPI = —————————         
             n 
           ∑         x weight
             i=1
'
$run = $cell.Characters(1, 13)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(14, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(15, 2)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(17, 12)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(29, 9)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Superscript = $true
$run = $cell.Characters(38, 13)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(51, 13)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(64, 3)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Superscript = $true
$run = $cell.Characters(67, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(68, 25)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(93, 13)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(106, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Subscript = $true
$run = $cell.Characters(107, 2)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(109, 12)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(121, 8)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Superscript = $true
$run = $cell.Characters(129, 10)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(139, 13)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(152, 3)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run.Font.Superscript = $true
$run = $cell.Characters(155, 1)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$ws.Range("C272").Value = 'Calculating the CPI for multiple items'
$cell = $ws.Range("D272")
$cell.Value = '         
Many but not all price indices are weighted averages using weights that sum to 1 or 100. Also the terms do not necessarily sum to 1 or 100.
Many but not all price indices are weighted averages using weights that sum to 1 or 100. 
Example: The prices of 85,000 items from 22,000 stores, and 35,000 rental units are added together and averaged. They are weighted this way: housing 41.4%; food and beverages 17.4%; transport 17.0%; medical care 6.9%; apparel 6.0%; entertainment 4.4%; other 6.9%. Taxes (43%) are not included in CPI computation.'
$run = $cell.Characters(1, 10)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(11, 140)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(151, 90)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$run = $cell.Characters(241, 312)
$run.Font.Size = 10
$run.Font.Name = "Helvetica Neue"
$ws.Range("E272").Value = 'Excerpt retrieved from:
en.wikipedia.org; 
title id 1062042105, Consumer Price Index 
This source was last edited on 25 Dec 2021, at 21:28 (UTC). 
Source above is available under the 
Creative Commons Attribution-ShareAlike License; additional terms may apply. By using this source, you agree to the Terms of Use and Privacy Policy. Wikipedia® is a registered trademark of the Wikimedia Foundation, Inc., a non-profit organization. 
'
$ws.Range("F272").Value = '.recapitulation'

# ---- Row 273 ----
$ws.Range("A273").Value = 'R053'
$ws.Range("B273").Value = 'HCPI'
$ws.Range("C273").Value = 'Harmonized Index of Consumer Prices (HICP)'
$ws.Range("D273").Value = 'By convention, weights are fractions or ratios summing to one, as percentages summing to 100 or as per mille numbers summing to 1000. 
On the European Union''s Harmonized Index of Consumer Prices (HICP), for example, each country computes some 80 prescribed sub-indices, their weighted average constituting the national HICP. The weights for these sub-indices will consist of the sum of the weights of a number of component lower level indices. The classification is according to use, developed in a national accounting context. This is not necessarily the kind of classification that is most appropriate for a consumer price index.'
$ws.Range("E273").Value = 'Excerpt retrieved from:
en.wikipedia.org; 
title id 1062042105, Consumer Price Index 
This source was last edited on 25 Dec 2021, at 21:28 (UTC). 
Source above is available under the 
Creative Commons Attribution-ShareAlike License; additional terms may apply. By using this source, you agree to the Terms of Use and Privacy Policy. Wikipedia® is a registered trademark of the Wikimedia Foundation, Inc., a non-profit organization.'
$ws.Range("F273").Value = '.recapitulation'

# ---- Apply the shared formatting (fill/border/alignment/style) used by
#      the rest of the table, copied from the last pre-existing data row ----
$ws.Range("A270:F270").Copy()
$ws.Range("A271:F271").PasteSpecial($xlPasteFormats)
$ws.Range("A270:F270").Copy()
$ws.Range("A272:F272").PasteSpecial($xlPasteFormats)
$ws.Range("A270:F270").Copy()
$ws.Range("A273:F273").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---- Row heights for the newly added rows ----
$ws.Rows("271").RowHeight = 146.7
$ws.Rows("272").RowHeight = 158.7
$ws.Rows("273").RowHeight = 170.7

# ---- Row-height side effects on pre-existing long-text rows, caused by
#      re-flowing the sheet after the edit (text wrap recalculation) ----
$ws.Rows("31").RowHeight = 1115.7
$ws.Rows("32").RowHeight = 634.2

